$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the marking/total figures (fix total marks error)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "78 / 112"
